$wb = $excel.ActiveWorkbook

# Update Metadata sheet timestamp
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 02:01 PM"

# Update Industry Analysis sheet 1 Year column (F2:F76)
$ws = $wb.Worksheets.Item("Industry Analysis")
$ws.Range("F2").Value = 21.0016
$ws.Range("F3").Value = -16.2396
$ws.Range("F4").Value = 27.1317
$ws.Range("F5").Value = -50.6494
$ws.Range("F6").Value = 53.2813
$ws.Range("F7").Value = -8.106199999999999
$ws.Range("F8").Value = -9.552099999999999
$ws.Range("F9").Value = 36.3756
$ws.Range("F10").Value = -6.1314
$ws.Range("F11").Value = 31.9081
$ws.Range("F12").Value = -18.4955
$ws.Range("F13").Value = 14.0155
$ws.Range("F14").Value = -36.0718
$ws.Range("F15").Value = -0.1622
$ws.Range("F16").Value = 0.1459
$ws.Range("F17").Value = -22.0012
$ws.Range("F18").Value = 1.0561
$ws.Range("F19").Value = -27.708
$ws.Range("F20").Value = 47.7309
$ws.Range("F21").Value = 12.0959
$ws.Range("F22").Value = 95.1491
$ws.Range("F23").Value = -50.2657
$ws.Range("F24").Value = -13.3427
$ws.Range("F25").Value = -9.9316
$ws.Range("F26").Value = 5.8244
$ws.Range("F27").Value = -32.7692
$ws.Range("F28").Value = -24.8224
$ws.Range("F29").Value = -18.4191
$ws.Range("F30").Value = 25.8569
$ws.Range("F31").Value = 58.4712
$ws.Range("F32").Value = -3.3862
$ws.Range("F33").Value = -6.3282
$ws.Range("F34").Value = 27.7203
$ws.Range("F35").Value = 4.4873
$ws.Range("F36").Value = -4.9458
$ws.Range("F37").Value = 3.6074
$ws.Range("F38").Value = -23.3973
$ws.Range("F39").Value = 8.7355
$ws.Range("F40").Value = -5.8541
$ws.Range("F41").Value = -8.3934
$ws.Range("F42").Value = 20.3818
$ws.Range("F43").Value = 14.3164
$ws.Range("F44").Value = -12.6846
$ws.Range("F45").Value = 28.4075
$ws.Range("F46").Value = -1.1135
$ws.Range("F47").Value = -37.1997
$ws.Range("F48").Value = -29.8569
$ws.Range("F49").Value = -27.5511
$ws.Range("F50").Value = -49.7478
$ws.Range("F51").Value = -51.8002
$ws.Range("F52").Value = -38.5254
$ws.Range("F53").Value = -12.4886
$ws.Range("F54").Value = -5.0725
$ws.Range("F55").Value = -17.7445
$ws.Range("F56").Value = -26.636
$ws.Range("F57").Value = -29.3361
$ws.Range("F58").Value = -11.9574
$ws.Range("F59").Value = -24.5687
$ws.Range("F60").Value = -12.3
$ws.Range("F61").Value = -10.9446
$ws.Range("F62").Value = -17.1229
$ws.Range("F63").Value = -9.5038
$ws.Range("F64").Value = 54.2749
$ws.Range("F65").Value = -43.4736
$ws.Range("F66").Value = 13.2687
$ws.Range("F67").Value = 12.7149
$ws.Range("F68").Value = 24.8057
$ws.Range("F69").Value = -17.0328
$ws.Range("F70").Value = -6.8927
$ws.Range("F71").Value = 13.6034
$ws.Range("F72").Value = 3.9995
$ws.Range("F73").Value = -16.226
$ws.Range("F74").Value = -16.2448
$ws.Range("F75").Value = 28.6924
$ws.Range("F76").Value = 48.9752
